$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version bump: 0.1.6 -> 0.1.7
$ws.Range("B3").Value = "0.1.7"

# Status: active -> draft
$ws.Range("B6").Value = "draft"

# Date update
$ws.Range("B8").Value = "2024-11-22T12:33:30-06:00"

# Contact row (row 10) now holds the publisher contact detail text
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

# Second Contact row (row 11) now holds the individual contact
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# Insert a new row for "Jurisdiction" (empty value) before the old "Description" row,
# pushing Description/Purpose/Copyright/Immutable down by one row. Copy the
# formatting from the row above so the new row matches the surrounding style.
$ws.Rows("12:12").Insert()
$ws.Range("A11:B11").Copy()
$ws.Range("A12:B12").PasteSpecial(-4122)
$ws.Range("A12").Value = "Jurisdiction"
$excel.CutCopyMode = 0
